# Applies the "require schema to be explicitly given in DB table
# definitions" edit to the DDS project workbook:
#   - Tables!B2:B7 "site"/"device"/"farm" -> "project.site"/"project.device"/"project.farm"
#   - Tables sheet becomes the active (selected) sheet/tab, with B7 selected
#   - ROOT sheet row 1 height reverts to the sheet default (was a custom 29pt row)

$wb = $excel.ActiveWorkbook

$tables = $wb.Worksheets.Item("Tables")
$root   = $wb.Worksheets.Item("ROOT")

# --- Qualify the DB table names referenced in column B of the Tables sheet ---
$tables.Range("B2").Value = "project.site"
$tables.Range("B3").Value = "project.device"
$tables.Range("B4").Value = "project.site"
$tables.Range("B5").Value = "project.site"
$tables.Range("B6").Value = "project.site"
$tables.Range("B7").Value = "project.farm"

# --- ROOT row 1 was a custom-height (29pt) header row; restore auto height ---
$root.Rows.Item(1).AutoFit()

# --- Make "Tables" the active sheet/tab, with B7 selected ---
$tables.Activate()
$tables.Range("B7").Select()
